$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.584.50'
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("D3").Value = '3.005.70'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.26'
$ws.Range("E5").Value = '  +2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.22'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '3.005.06'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.92'
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("E12").Value = '  +3.66%  '
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("E15").Value = '  +2.36%  '
$ws.Range("D16").Value = '3.501.96'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.03'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '61.553.27'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = '3.005.26'
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '454.04'
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.04'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.38'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.92'
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -4.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.68'
$ws.Range("E26").Value = '  +6.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.99'
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.68'
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.25'
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.50'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.108'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").Value = '0.0₃0839'
$ws.Range("E35").Value = '  +5.52%  '
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.23'
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("E39").Value = '  -2.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.37'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  +9.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.92'
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '401.36'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.98'
$ws.Range("E44").Value = '  +5.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0354'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.271'
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("D47").Value = '2.720.36'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.60'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("E50").Value = '  -0.47%  '
